$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextCell 'D2' '57.620.26'
Set-TextCell 'E2' '  +0.27%  '
Set-TextCell 'D3' '3.119.97'
Set-TextCell 'E3' '  +0.45%  '
Set-TextCell 'E4' '  +0.00%  '
Set-TextCell 'D5' '531.32'
Set-TextCell 'E5' '  +1.30%  '
Set-TextCell 'D6' '138.22'
Set-TextCell 'E6' '  +1.33%  '
Set-TextCell 'E7' '  +0.02%  '
Set-TextCell 'D8' '3.118.29'
Set-TextCell 'E8' '  +0.44%  '
Set-TextCell 'D9' '0.471'
Set-TextCell 'E9' '  +5.53%  '
Set-TextCell 'D10' '7.31'
Set-TextCell 'E10' '  +0.23%  '
Set-TextCell 'E11' '  +0.52%  '
Set-TextCell 'E13' '  +1.46%  '
Set-TextCell 'D14' '3.653.88'
Set-TextCell 'E14' '  +0.27%  '
Set-TextCell 'D15' '25.68'
Set-TextCell 'E15' '  +2.05%  '
Set-TextCell 'E16' '  +1.17%  '
Set-TextCell 'D17' '57.721.99'
Set-TextCell 'E17' '  +0.31%  '
Set-TextCell 'D18' '3.116.90'
Set-TextCell 'E18' '  +0.38%  '
Set-TextCell 'D19' '6.04'
Set-TextCell 'E19' '  +2.11%  '
Set-TextCell 'E20' '  +2.37%  '
Set-TextCell 'D21' '8.07'
Set-TextCell 'E21' '  +2.75%  '
Set-TextCell 'D22' '361.68'
Set-TextCell 'E22' '  +4.69%  '
Set-TextCell 'E23' '  -0.01%  '
Set-TextCell 'D24' '68.98'
Set-TextCell 'E24' '  +1.99%  '
Set-TextCell 'E25' '  +0.89%  '
Set-TextCell 'E26' '  +0.19%  '
Set-TextCell 'E27' '  -0.13%  '
Set-TextCell 'D28' '0.0₃0863'
Set-TextCell 'E28' '  -2.78%  '
Set-TextCell 'E29' '  -1.28%  '
Set-TextCell 'B30' 'PancakeSwap'
Set-TextCell 'C30' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D30' '1.87'
Set-TextCell 'E30' '  +0.17%  '
Set-TextCell 'B31' 'RenderToken'
Set-TextCell 'C31' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D31' '6.07'
Set-TextCell 'E31' '  +0.32%  '
Set-TextCell 'E32' '  +2.59%  '
Set-TextCell 'D33' '5.13'
Set-TextCell 'E33' '  +4.50%  '
Set-TextCell 'E34' '  -0.34%  '
Set-TextCell 'D35' '159.31'
Set-TextCell 'E35' '  +0.73%  '
Set-TextCell 'E36' '  +0.32%  '
Set-TextCell 'E37' '  +4.85%  '
Set-TextCell 'D38' '25.42'
Set-TextCell 'E38' '  -1.43%  '
Set-TextCell 'E39' '  +4.20%  '
Set-TextCell 'E40' '  +1.64%  '
Set-TextCell 'D41' '2.540.81'
Set-TextCell 'E41' '  +7.79%  '
Set-TextCell 'D42' '4.02'
Set-TextCell 'E42' '  -1.84%  '
Set-TextCell 'D43' '0.698'
Set-TextCell 'E43' '  -0.13%  '
Set-TextCell 'D44' '37.81'
Set-TextCell 'D45' '0.0269'
Set-TextCell 'E45' '  +1.23%  '
Set-TextCell 'D46' '0.999'
Set-TextCell 'E46' '  -0.04%  '
Set-TextCell 'D47' '0.977'
Set-TextCell 'E47' '  +1.26%  '
Set-TextCell 'D48' '6.08'
Set-TextCell 'E48' '  +1.98%  '
Set-TextCell 'D49' '19.69'
Set-TextCell 'E49' '  -0.41%  '
Set-TextCell 'D50' '0.742'
Set-TextCell 'E50' '  -2.38%  '
Set-TextCell 'D51' '0.0913'
Set-TextCell 'E51' '  +2.91%  '
